$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values (rows 2-3) ---
$ws.Range("B2").Value = 100
$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 1

# --- Apply the same cell formatting (style) used by the existing data rows (2-3) ---
# to the previously-empty data rows 4-7, columns A:I.
$ws.Range("B2").Copy()
$ws.Range("A4:I7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the row heights used by the other data rows.
$ws.Rows.Item(4).RowHeight = 16.5
$ws.Rows.Item(5).RowHeight = 16.5
$ws.Rows.Item(6).RowHeight = 16.5
$ws.Rows.Item(7).RowHeight = 16.5

# --- Fill in the new scenario rows ---
# Row 4 : id=2
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 200
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 100
$ws.Range("G4").Value = 0.5
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 0.8

# Row 5 : id=3
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 100
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 200
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = 0.5
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 0.9

# Row 6 : id=4
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 100
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 200
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 100
$ws.Range("G6").Value = 0.3
$ws.Range("H6").Value = 100
$ws.Range("I6").Value = 0.9

# Row 7 : id=5
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 100
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 200
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 100
$ws.Range("G7").Value = 0.7
$ws.Range("H7").Value = 100
$ws.Range("I7").Value = 0.9

# --- Cosmetic: update the active selection like in the saved file ---
$null = $ws.Range("G8").Select()
